$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 11:34:59"
$ws1.Range("A3").Value = "Total filas: 132"
$ws1.Range("A102").Value = "11:34:59"
$ws1.Range("B102").Value = "11:44"
$ws1.Range("C102").Value = "15X38_ABASTO"
$ws1.Range("D102").Value = 10
$ws1.Range("E102").Value = "LP1912"
$ws1.Range("A103").Value = "10:05:51"
$ws1.Range("B103").Value = "11:45"
$ws1.Range("C103").Value = "15X38_ABASTO"
$ws1.Range("D103").Value = 100
$ws1.Range("E103").Value = "LP1912"
$ws1.Range("A104").Value = "11:11:33"
$ws1.Range("B104").Value = "11:51"
$ws1.Range("C104").Value = "23_HERNANDEZ"
$ws1.Range("D104").Value = 40
$ws1.Range("E104").Value = "LP1912"
$ws1.Range("A105").Value = "10:05:51"
$ws1.Range("B105").Value = "11:52"
$ws1.Range("C105").Value = "225_GOMEZ"
$ws1.Range("D105").Value = 107
$ws1.Range("E105").Value = "LP1912"
$ws1.Range("A106").Value = "10:50:41"
$ws1.Range("B106").Value = "11:53"
$ws1.Range("C106").Value = "225_GOMEZ"
$ws1.Range("D106").Value = 63
$ws1.Range("E106").Value = "LP1912"
$ws1.Range("A107").Value = "10:37:52"
$ws1.Range("B107").Value = "11:53"
$ws1.Range("C107").Value = "23_HERNANDEZ"
$ws1.Range("D107").Value = 76
$ws1.Range("E107").Value = "LP1912"
$ws1.Range("A108").Value = "10:50:41"
$ws1.Range("B108").Value = "11:54"
$ws1.Range("C108").Value = "23_HERNANDEZ"
$ws1.Range("D108").Value = 64
$ws1.Range("E108").Value = "LP1912"
$ws1.Range("A109").Value = "11:34:59"
$ws1.Range("B109").Value = "11:57"
$ws1.Range("C109").Value = "17_ROMERO"
$ws1.Range("D109").Value = 23
$ws1.Range("E109").Value = "LP1912"
$ws1.Range("A110").Value = "10:05:51"
$ws1.Range("B110").Value = "11:58"
$ws1.Range("C110").Value = "17_ROMERO"
$ws1.Range("D110").Value = 113
$ws1.Range("E110").Value = "LP1912"
$ws1.Range("A111").Value = "10:37:52"
$ws1.Range("B111").Value = "12:05"
$ws1.Range("C111").Value = "11_ETCHEVERRY"
$ws1.Range("D111").Value = 88
$ws1.Range("E111").Value = "LP1912"
$ws1.Range("A112").Value = "11:34:59"
$ws1.Range("B112").Value = "12:09"
$ws1.Range("C112").Value = "15_ABASTO"
$ws1.Range("D112").Value = 35
$ws1.Range("E112").Value = "LP1912"
$ws1.Range("A113").Value = "11:34:59"
$ws1.Range("B113").Value = "12:09"
$ws1.Range("C113").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D113").Value = 35
$ws1.Range("E113").Value = "LP1912"
$ws1.Range("A114").Value = "10:37:52"
$ws1.Range("B114").Value = "12:10"
$ws1.Range("C114").Value = "15_ABASTO"
$ws1.Range("D114").Value = 93
$ws1.Range("E114").Value = "LP1912"
$ws1.Range("A115").Value = "10:37:52"
$ws1.Range("B115").Value = "12:10"
$ws1.Range("C115").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D115").Value = 93
$ws1.Range("E115").Value = "LP1912"
$ws1.Range("A116").Value = "10:37:52"
$ws1.Range("B116").Value = "12:16"
$ws1.Range("C116").Value = "10_OLMOS"
$ws1.Range("D116").Value = 99
$ws1.Range("E116").Value = "LP1912"
$ws1.Range("A117").Value = "11:11:33"
$ws1.Range("B117").Value = "12:17"
$ws1.Range("C117").Value = "10_OLMOS"
$ws1.Range("D117").Value = 66
$ws1.Range("E117").Value = "LP1912"
$ws1.Range("A118").Value = "10:37:52"
$ws1.Range("B118").Value = "12:21"
$ws1.Range("C118").Value = "215C_EL PATO"
$ws1.Range("D118").Value = 104
$ws1.Range("E118").Value = "LP1912"
$ws1.Range("A119").Value = "11:11:33"
$ws1.Range("B119").Value = "12:22"
$ws1.Range("C119").Value = "215C_EL PATO"
$ws1.Range("D119").Value = 71
$ws1.Range("E119").Value = "LP1912"
$ws1.Range("A120").Value = "10:37:52"
$ws1.Range("B120").Value = "12:32"
$ws1.Range("C120").Value = "14_ABASTO"
$ws1.Range("D120").Value = 115
$ws1.Range("E120").Value = "LP1912"
$ws1.Range("A121").Value = "11:34:59"
$ws1.Range("B121").Value = "12:33"
$ws1.Range("C121").Value = "15_ABASTO"
$ws1.Range("D121").Value = 59
$ws1.Range("E121").Value = "LP1912"
$ws1.Range("A122").Value = "10:37:52"
$ws1.Range("B122").Value = "12:34"
$ws1.Range("C122").Value = "15_ABASTO"
$ws1.Range("D122").Value = 117
$ws1.Range("E122").Value = "LP1912"
$ws1.Range("A123").Value = "11:11:33"
$ws1.Range("B123").Value = "12:35"
$ws1.Range("C123").Value = "23_HERNANDEZ"
$ws1.Range("D123").Value = 84
$ws1.Range("E123").Value = "LP1912"
$ws1.Range("A124").Value = "11:34:59"
$ws1.Range("B124").Value = "12:35"
$ws1.Range("C124").Value = "27_EL RETIRO"
$ws1.Range("D124").Value = 61
$ws1.Range("E124").Value = "LP1912"
$ws1.Range("A125").Value = "10:50:41"
$ws1.Range("B125").Value = "12:36"
$ws1.Range("C125").Value = "27_EL RETIRO"
$ws1.Range("D125").Value = 106
$ws1.Range("E125").Value = "LP1912"
$ws1.Range("A126").Value = "11:34:59"
$ws1.Range("B126").Value = "12:36"
$ws1.Range("C126").Value = "23_HERNANDEZ"
$ws1.Range("D126").Value = 62
$ws1.Range("E126").Value = "LP1912"
$ws1.Range("A127").Value = "11:34:59"
$ws1.Range("B127").Value = "12:47"
$ws1.Range("C127").Value = "15X38_ABASTO"
$ws1.Range("D127").Value = 73
$ws1.Range("E127").Value = "LP1912"
$ws1.Range("A128").Value = "11:34:59"
$ws1.Range("B128").Value = "12:47"
$ws1.Range("C128").Value = "16_SANTA ANA"
$ws1.Range("D128").Value = 73
$ws1.Range("E128").Value = "LP1912"
$ws1.Range("A129").Value = "11:34:59"
$ws1.Range("B129").Value = "12:47"
$ws1.Range("C129").Value = "14_ABASTO"
$ws1.Range("D129").Value = 73
$ws1.Range("E129").Value = "LP1912"
$ws1.Range("A130").Value = "11:11:33"
$ws1.Range("B130").Value = "12:48"
$ws1.Range("C130").Value = "15X38_ABASTO"
$ws1.Range("D130").Value = 97
$ws1.Range("E130").Value = "LP1912"
$ws1.Range("A131").Value = "10:50:41"
$ws1.Range("B131").Value = "12:48"
$ws1.Range("C131").Value = "16_SANTA ANA"
$ws1.Range("D131").Value = 118
$ws1.Range("E131").Value = "LP1912"
$ws1.Range("A132").Value = "11:11:33"
$ws1.Range("B132").Value = "13:02"
$ws1.Range("C132").Value = "11_ETCHEVERRY"
$ws1.Range("D132").Value = 111
$ws1.Range("E132").Value = "LP1912"
$ws1.Range("A133").Value = "11:34:59"
$ws1.Range("B133").Value = "13:03"
$ws1.Range("C133").Value = "215C_EL PATO"
$ws1.Range("D133").Value = 89
$ws1.Range("E133").Value = "LP1912"
$ws1.Range("A134").Value = "11:34:59"
$ws1.Range("B134").Value = "13:12"
$ws1.Range("C134").Value = "16_SANTA ANA"
$ws1.Range("D134").Value = 98
$ws1.Range("E134").Value = "LP1912"
$ws1.Range("A135").Value = "11:34:59"
$ws1.Range("B135").Value = "13:16"
$ws1.Range("C135").Value = "10_OLMOS"
$ws1.Range("D135").Value = 102
$ws1.Range("E135").Value = "LP1912"
$ws1.Range("A136").Value = "11:34:59"
$ws1.Range("B136").Value = "13:24"
$ws1.Range("C136").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D136").Value = 110
$ws1.Range("E136").Value = "LP1912"
$ws1.Range("A137").Value = "11:34:59"
$ws1.Range("B137").Value = "13:32"
$ws1.Range("C137").Value = "215A_EL PATO"
$ws1.Range("D137").Value = 118
$ws1.Range("E137").Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 11:34:59"
$ws2.Range("A3").Value = "Total filas: 24"
$ws2.Range("A28").Value = "11:34:59"
$ws2.Range("B28").Value = "13:03"
$ws2.Range("C28").Value = "215C_EL PATO"
$ws2.Range("D28").Value = 89
$ws2.Range("E28").Value = "LP1912"
$ws2.Range("A29").Value = "11:34:59"
$ws2.Range("B29").Value = "13:32"
$ws2.Range("C29").Value = "215A_EL PATO"
$ws2.Range("D29").Value = 118
$ws2.Range("E29").Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 11:34:59"
$ws3.Range("A3").Value = "Total filas: 19"
$ws3.Range("A23").Value = "11:34:59"
$ws3.Range("B23").Value = "13:11"
$ws3.Range("C23").Value = "215C_LA PLATA"
$ws3.Range("D23").Value = 97
$ws3.Range("E23").Value = "L6203"
$ws3.Range("A24").Value = "11:34:59"
$ws3.Range("B24").Value = "13:20"
$ws3.Range("C24").Value = "215B_LP-P MOR-1 Y 57"
$ws3.Range("D24").Value = 106
$ws3.Range("E24").Value = "L6173"
